$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
# row 264
$ws.Cells.Item(264, 1).Value = 45440.59375
$ws.Cells.Item(264, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(264, 2).Value = '13-06-2024 10:15:00'
$ws.Cells.Item(264, 3).Value = 'hour'
$ws.Cells.Item(264, 4).Value = 'RAIN.NS'
$ws.Cells.Item(264, 5).Value = 45408.46875
$ws.Cells.Item(264, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(264, 6).Value = 184.5
$ws.Cells.Item(264, 7).Value = 45434.42708333334
$ws.Cells.Item(264, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(264, 8).Value = 174.75
$ws.Cells.Item(264, 9).Value = 45439.42708333334
$ws.Cells.Item(264, 9).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(264, 10).Value = 173.6499938964844
$ws.Cells.Item(264, 11).Value = 'High'
$ws.Cells.Item(264, 12).Value = '13/06/2024 04:46:56'

# row 265
$ws.Cells.Item(265, 1).Value = 45411.55208333334
$ws.Cells.Item(265, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(265, 2).Value = '13-06-2024 09:15:00'
$ws.Cells.Item(265, 3).Value = 'hour'
$ws.Cells.Item(265, 4).Value = 'DEVYANI.NS'
$ws.Cells.Item(265, 5).Value = 45391.46875
$ws.Cells.Item(265, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(265, 6).Value = 167.8000030517578
$ws.Cells.Item(265, 7).Value = 45394.38541666666
$ws.Cells.Item(265, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(265, 8).Value = 168
$ws.Cells.Item(265, 9).Value = 45408.38541666666
$ws.Cells.Item(265, 9).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(265, 10).Value = 170.6000061035156
$ws.Cells.Item(265, 11).Value = 'High'
$ws.Cells.Item(265, 12).Value = '13/06/2024 04:46:56'


# Sheet 2
$ws = $wb.Worksheets.Item(2)
# row 86
$ws.Cells.Item(86, 1).Value = 45447.55208333334
$ws.Cells.Item(86, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(86, 2).Value = '13-06-2024 10:15:00'
$ws.Cells.Item(86, 3).Value = 'hour'
$ws.Cells.Item(86, 4).Value = 'KOTAKBANK.NS'
$ws.Cells.Item(86, 5).Value = 45439.51041666666
$ws.Cells.Item(86, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(86, 6).Value = 1722
$ws.Cells.Item(86, 7).Value = 45446.38541666666
$ws.Cells.Item(86, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(86, 8).Value = 1726.449951171875
$ws.Cells.Item(86, 9).Value = 'High'
$ws.Cells.Item(86, 10).Value = '13/06/2024 04:46:56'

# row 87
$ws.Cells.Item(87, 1).Value = 45429.55208333334
$ws.Cells.Item(87, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(87, 2).Value = '13-06-2024 09:15:00'
$ws.Cells.Item(87, 3).Value = 'hour'
$ws.Cells.Item(87, 4).Value = 'KSOLVES.NS'
$ws.Cells.Item(87, 5).Value = 45426.55208333334
$ws.Cells.Item(87, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(87, 6).Value = 1135
$ws.Cells.Item(87, 7).Value = 45428.38541666666
$ws.Cells.Item(87, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(87, 8).Value = 1135
$ws.Cells.Item(87, 9).Value = 'High'
$ws.Cells.Item(87, 10).Value = '13/06/2024 04:46:56'

# row 88
$ws.Cells.Item(88, 1).Value = 45433.59375
$ws.Cells.Item(88, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(88, 2).Value = '13-06-2024 09:15:00'
$ws.Cells.Item(88, 3).Value = 'hour'
$ws.Cells.Item(88, 4).Value = 'WALCHANNAG.NS'
$ws.Cells.Item(88, 5).Value = 45429.38541666666
$ws.Cells.Item(88, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(88, 6).Value = 219.0500030517578
$ws.Cells.Item(88, 7).Value = 45429.42708333334
$ws.Cells.Item(88, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(88, 8).Value = 219.0500030517578
$ws.Cells.Item(88, 9).Value = 'High'
$ws.Cells.Item(88, 10).Value = '13/06/2024 04:46:56'

# row 89
$ws.Cells.Item(89, 1).Value = 45433.63541666666
$ws.Cells.Item(89, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(89, 2).Value = '13-06-2024 09:15:00'
$ws.Cells.Item(89, 3).Value = 'hour'
$ws.Cells.Item(89, 4).Value = 'WALCHANNAG.NS'
$ws.Cells.Item(89, 5).Value = 45429.38541666666
$ws.Cells.Item(89, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(89, 6).Value = 219.0500030517578
$ws.Cells.Item(89, 7).Value = 45429.46875
$ws.Cells.Item(89, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(89, 8).Value = 219.0500030517578
$ws.Cells.Item(89, 9).Value = 'High'
$ws.Cells.Item(89, 10).Value = '13/06/2024 04:46:56'


# Sheet 3
$ws = $wb.Worksheets.Item(3)
# row 840
$ws.Cells.Item(840, 1).Value = 'CROWN.NS'
$ws.Cells.Item(840, 2).Value = 45446.38541666666
$ws.Cells.Item(840, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(840, 3).Value = 244.9499969482422
$ws.Cells.Item(840, 4).Value = 244.9499969482422
$ws.Cells.Item(840, 5).Value = 244.9499969482422
$ws.Cells.Item(840, 6).Value = 'High'
$ws.Cells.Item(840, 7).Value = 244.9499969482422
$ws.Cells.Item(840, 8).Value = 'hour'
$ws.Cells.Item(840, 9).Value = '13-06-2024 09:15:00'
$ws.Cells.Item(840, 10).Value = 245.0899963378906
$ws.Cells.Item(840, 11).Value = 240.2899932861328
$ws.Cells.Item(840, 12).Value = '13/06/2024 04:46:56'

# row 841
$ws.Cells.Item(841, 1).Value = 'ZODIAC.NS'
$ws.Cells.Item(841, 2).Value = 45448.51041666666
$ws.Cells.Item(841, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(841, 3).Value = 611.8499755859375
$ws.Cells.Item(841, 4).Value = 611.8499755859375
$ws.Cells.Item(841, 5).Value = 611.8499755859375
$ws.Cells.Item(841, 6).Value = 'Low'
$ws.Cells.Item(841, 7).Value = 611.8499755859375
$ws.Cells.Item(841, 8).Value = 'hour'
$ws.Cells.Item(841, 9).Value = '13-06-2024 09:15:00'
$ws.Cells.Item(841, 10).Value = 611.0999755859375
$ws.Cells.Item(841, 11).Value = 623.5499877929688
$ws.Cells.Item(841, 12).Value = '13/06/2024 04:46:56'

# row 842
$ws.Cells.Item(842, 1).Value = 'ZODIAC.NS'
$ws.Cells.Item(842, 2).Value = 45448.55208333334
$ws.Cells.Item(842, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(842, 3).Value = 611.8499755859375
$ws.Cells.Item(842, 4).Value = 611.8499755859375
$ws.Cells.Item(842, 5).Value = 611.8499755859375
$ws.Cells.Item(842, 6).Value = 'Low'
$ws.Cells.Item(842, 7).Value = 611.8499755859375
$ws.Cells.Item(842, 8).Value = 'hour'
$ws.Cells.Item(842, 9).Value = '13-06-2024 09:15:00'
$ws.Cells.Item(842, 10).Value = 611.0999755859375
$ws.Cells.Item(842, 11).Value = 623.5499877929688
$ws.Cells.Item(842, 12).Value = '13/06/2024 04:46:56'

# row 843
$ws.Cells.Item(843, 1).Value = 'ZODIAC.NS'
$ws.Cells.Item(843, 2).Value = 45448.59375
$ws.Cells.Item(843, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(843, 3).Value = 611.8499755859375
$ws.Cells.Item(843, 4).Value = 611.8499755859375
$ws.Cells.Item(843, 5).Value = 611.8499755859375
$ws.Cells.Item(843, 6).Value = 'Low'
$ws.Cells.Item(843, 7).Value = 611.8499755859375
$ws.Cells.Item(843, 8).Value = 'hour'
$ws.Cells.Item(843, 9).Value = '13-06-2024 09:15:00'
$ws.Cells.Item(843, 10).Value = 611.0999755859375
$ws.Cells.Item(843, 11).Value = 623.5499877929688
$ws.Cells.Item(843, 12).Value = '13/06/2024 04:46:56'

# row 844
$ws.Cells.Item(844, 1).Value = 'SPELS.BO'
$ws.Cells.Item(844, 2).Value = 45434.42708333334
$ws.Cells.Item(844, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(844, 3).Value = 141.6999969482422
$ws.Cells.Item(844, 4).Value = 141.6999969482422
$ws.Cells.Item(844, 5).Value = 141.6999969482422
$ws.Cells.Item(844, 6).Value = 'High'
$ws.Cells.Item(844, 7).Value = 141.6999969482422
$ws.Cells.Item(844, 8).Value = 'hour'
$ws.Cells.Item(844, 9).Value = '13-06-2024 09:15:00'
$ws.Cells.Item(844, 10).Value = 142.5
$ws.Cells.Item(844, 11).Value = 138.6499938964844
$ws.Cells.Item(844, 12).Value = '13/06/2024 04:46:56'

# row 845
$ws.Cells.Item(845, 1).Value = 'SPELS.BO'
$ws.Cells.Item(845, 2).Value = 45434.46875
$ws.Cells.Item(845, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(845, 3).Value = 141.6999969482422
$ws.Cells.Item(845, 4).Value = 141.6999969482422
$ws.Cells.Item(845, 5).Value = 141.6999969482422
$ws.Cells.Item(845, 6).Value = 'High'
$ws.Cells.Item(845, 7).Value = 141.6999969482422
$ws.Cells.Item(845, 8).Value = 'hour'
$ws.Cells.Item(845, 9).Value = '13-06-2024 09:15:00'
$ws.Cells.Item(845, 10).Value = 142.5
$ws.Cells.Item(845, 11).Value = 138.6499938964844
$ws.Cells.Item(845, 12).Value = '13/06/2024 04:46:56'

# row 846
$ws.Cells.Item(846, 1).Value = 'CENTRALBK.BO'
$ws.Cells.Item(846, 2).Value = 45449.46875
$ws.Cells.Item(846, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(846, 3).Value = 66.45999908447266
$ws.Cells.Item(846, 4).Value = 64.69999694824219
$ws.Cells.Item(846, 5).Value = 64.93000030517578
$ws.Cells.Item(846, 6).Value = 'High'
$ws.Cells.Item(846, 7).Value = 66.45999908447266
$ws.Cells.Item(846, 8).Value = 'hour'
$ws.Cells.Item(846, 9).Value = '13-06-2024 09:15:00'
$ws.Cells.Item(846, 10).Value = 67.09999847412109
$ws.Cells.Item(846, 11).Value = 66.43000030517578
$ws.Cells.Item(846, 12).Value = '13/06/2024 04:46:56'

# row 847
$ws.Cells.Item(847, 1).Value = 'WINSOMTX.BO'
$ws.Cells.Item(847, 2).Value = 45446.38541666666
$ws.Cells.Item(847, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(847, 3).Value = 80.90000152587891
$ws.Cells.Item(847, 4).Value = 79
$ws.Cells.Item(847, 5).Value = 79
$ws.Cells.Item(847, 6).Value = 'High'
$ws.Cells.Item(847, 7).Value = 80.90000152587891
$ws.Cells.Item(847, 8).Value = 'hour'
$ws.Cells.Item(847, 9).Value = '13-06-2024 09:15:00'
$ws.Cells.Item(847, 10).Value = 81.5
$ws.Cells.Item(847, 11).Value = 80.09999847412109
$ws.Cells.Item(847, 12).Value = '13/06/2024 04:46:56'

# row 848
$ws.Cells.Item(848, 1).Value = 'HTMEDIA.NS'
$ws.Cells.Item(848, 2).Value = 45436.38541666666
$ws.Cells.Item(848, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(848, 3).Value = 27.10000038146973
$ws.Cells.Item(848, 4).Value = 26.60000038146973
$ws.Cells.Item(848, 5).Value = 26.85000038146973
$ws.Cells.Item(848, 6).Value = 'High'
$ws.Cells.Item(848, 7).Value = 27.10000038146973
$ws.Cells.Item(848, 8).Value = 'hour'
$ws.Cells.Item(848, 9).Value = '13-06-2024 09:15:00'
$ws.Cells.Item(848, 10).Value = 27.11000061035156
$ws.Cells.Item(848, 11).Value = 26.84000015258789
$ws.Cells.Item(848, 12).Value = '13/06/2024 04:46:56'

# row 849
$ws.Cells.Item(849, 1).Value = 'HTMEDIA.NS'
$ws.Cells.Item(849, 2).Value = 45454.42708333334
$ws.Cells.Item(849, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(849, 3).Value = 26.94000053405762
$ws.Cells.Item(849, 4).Value = 26.54999923706055
$ws.Cells.Item(849, 5).Value = 26.57999992370605
$ws.Cells.Item(849, 6).Value = 'High'
$ws.Cells.Item(849, 7).Value = 26.94000053405762
$ws.Cells.Item(849, 8).Value = 'hour'
$ws.Cells.Item(849, 9).Value = '13-06-2024 09:15:00'
$ws.Cells.Item(849, 10).Value = 27.11000061035156
$ws.Cells.Item(849, 11).Value = 26.84000015258789
$ws.Cells.Item(849, 12).Value = '13/06/2024 04:46:56'

# row 850
$ws.Cells.Item(850, 1).Value = 'KPEL.BO'
$ws.Cells.Item(850, 2).Value = 45446.38541666666
$ws.Cells.Item(850, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(850, 3).Value = 436.2000122070312
$ws.Cells.Item(850, 4).Value = 423.7999877929688
$ws.Cells.Item(850, 5).Value = 428.8999938964844
$ws.Cells.Item(850, 6).Value = 'High'
$ws.Cells.Item(850, 7).Value = 436.2000122070312
$ws.Cells.Item(850, 8).Value = 'hour'
$ws.Cells.Item(850, 9).Value = '13-06-2024 09:15:00'
$ws.Cells.Item(850, 10).Value = 437.3999938964844
$ws.Cells.Item(850, 11).Value = 434
$ws.Cells.Item(850, 12).Value = '13/06/2024 04:46:56'

# row 851
$ws.Cells.Item(851, 1).Value = 'ADFFOODS.NS'
$ws.Cells.Item(851, 2).Value = 45446.38541666666
$ws.Cells.Item(851, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(851, 3).Value = 227.25
$ws.Cells.Item(851, 4).Value = 215.3999938964844
$ws.Cells.Item(851, 5).Value = 217.75
$ws.Cells.Item(851, 6).Value = 'High'
$ws.Cells.Item(851, 7).Value = 227.25
$ws.Cells.Item(851, 8).Value = 'hour'
$ws.Cells.Item(851, 9).Value = '13-06-2024 09:15:00'
$ws.Cells.Item(851, 10).Value = 227.5099945068359
$ws.Cells.Item(851, 11).Value = 224.0099945068359
$ws.Cells.Item(851, 12).Value = '13/06/2024 04:46:56'

